# Generate Report for Handoff
# The 0ef9996f file moved from "In Translation" to "Ready for handoff" (with a new
# handoff timestamp). As a result, the two files swap places (rows 2 and 3) in each
# of the report sheets, and the hyperlinks' display text moves with the swapped
# values while their underlying relationship ids stay tied to the worksheet rows.

$wb = $excel.ActiveWorkbook

function Set-LinkDisplay($ws, $cellAddress, $text) {
    foreach ($hl in $ws.Hyperlinks) {
        if ($hl.Range.Address() -eq $cellAddress) {
            $hl.TextToDisplay = $text
        }
    }
}

# ---------------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A2").Value = "77d4fc19-41dc-4c65-a058-103f8a9d2c26.md"
$wsOverview.Range("B2").Value = "In Translation"
$wsOverview.Range("C2").Value = "In Translation"
$wsOverview.Range("D2").Value = "2016-16-20 12:16:12"

$wsOverview.Range("A3").Value = "0ef9996f-6e2a-4f59-92ae-f5e5ce48e193.md"
$wsOverview.Range("B3").Value = "Ready for handoff"
$wsOverview.Range("C3").Value = "Ready for handoff"
$wsOverview.Range("D3").Value = "2016-16-20 12:16:56"

Set-LinkDisplay $wsOverview '$A$2' "77d4fc19-41dc-4c65-a058-103f8a9d2c26.md"
Set-LinkDisplay $wsOverview '$A$3' "0ef9996f-6e2a-4f59-92ae-f5e5ce48e193.md"

# ---------------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("A2").Value = "77d4fc19-41dc-4c65-a058-103f8a9d2c26.md"
$wsZh.Range("B2").Value = ".md"
$wsZh.Range("C2").Value = "In Translation"
$wsZh.Range("D2").Value = "77d4fc19-41dc-4c65-a058-103f8a9d2c26.42bb164d26d94036fca13ac4d4563e1f448c82a1.zh-cn.xlf"
$wsZh.Range("E2").Value = "2016-03-20 12:16:09"
$wsZh.Range("H2").Value = "0001-01-01 00:00:00"
$wsZh.Range("I2").Value = "Include"

$wsZh.Range("A3").Value = "0ef9996f-6e2a-4f59-92ae-f5e5ce48e193.md"
$wsZh.Range("B3").Value = ".md"
$wsZh.Range("C3").Value = "Ready for handoff"
$wsZh.Range("D3").Value = "0ef9996f-6e2a-4f59-92ae-f5e5ce48e193.bc7d1377f5ea65fc4dd773bfd88222caa9bb8e8b.zh-cn.xlf"
$wsZh.Range("E3").Value = "2016-03-20 12:16:53"
$wsZh.Range("H3").Value = "0001-01-01 00:00:00"
$wsZh.Range("I3").Value = "Include"

Set-LinkDisplay $wsZh '$A$2' "77d4fc19-41dc-4c65-a058-103f8a9d2c26.md"
Set-LinkDisplay $wsZh '$B$2' ".md"
Set-LinkDisplay $wsZh '$D$2' "77d4fc19-41dc-4c65-a058-103f8a9d2c26.42bb164d26d94036fca13ac4d4563e1f448c82a1.zh-cn.xlf"
Set-LinkDisplay $wsZh '$A$3' "0ef9996f-6e2a-4f59-92ae-f5e5ce48e193.md"
Set-LinkDisplay $wsZh '$B$3' ".md"
Set-LinkDisplay $wsZh '$D$3' "0ef9996f-6e2a-4f59-92ae-f5e5ce48e193.bc7d1377f5ea65fc4dd773bfd88222caa9bb8e8b.zh-cn.xlf"

# ---------------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("A2").Value = "77d4fc19-41dc-4c65-a058-103f8a9d2c26.md"
$wsDe.Range("B2").Value = ".md"
$wsDe.Range("C2").Value = "In Translation"
$wsDe.Range("D2").Value = "77d4fc19-41dc-4c65-a058-103f8a9d2c26.42bb164d26d94036fca13ac4d4563e1f448c82a1.de-de.xlf"
$wsDe.Range("E2").Value = "2016-03-20 12:16:12"
$wsDe.Range("H2").Value = "0001-01-01 00:00:00"
$wsDe.Range("I2").Value = "Include"

$wsDe.Range("A3").Value = "0ef9996f-6e2a-4f59-92ae-f5e5ce48e193.md"
$wsDe.Range("B3").Value = ".md"
$wsDe.Range("C3").Value = "Ready for handoff"
$wsDe.Range("D3").Value = "0ef9996f-6e2a-4f59-92ae-f5e5ce48e193.bc7d1377f5ea65fc4dd773bfd88222caa9bb8e8b.de-de.xlf"
$wsDe.Range("E3").Value = "2016-03-20 12:16:56"
$wsDe.Range("H3").Value = "0001-01-01 00:00:00"
$wsDe.Range("I3").Value = "Include"

Set-LinkDisplay $wsDe '$A$2' "77d4fc19-41dc-4c65-a058-103f8a9d2c26.md"
Set-LinkDisplay $wsDe '$B$2' ".md"
Set-LinkDisplay $wsDe '$D$2' "77d4fc19-41dc-4c65-a058-103f8a9d2c26.42bb164d26d94036fca13ac4d4563e1f448c82a1.de-de.xlf"
Set-LinkDisplay $wsDe '$A$3' "0ef9996f-6e2a-4f59-92ae-f5e5ce48e193.md"
Set-LinkDisplay $wsDe '$B$3' ".md"
Set-LinkDisplay $wsDe '$D$3' "0ef9996f-6e2a-4f59-92ae-f5e5ce48e193.bc7d1377f5ea65fc4dd773bfd88222caa9bb8e8b.de-de.xlf"
